# Append newly fetched ticker rows to the end of the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTickers = @("IMX-USD", "MNT-USD", "PEPE-USD", "GRT-USD", "TAO-USD")

# Find the next empty row after the current last used row (row 468 -> starts at 469).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newTickers.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newTickers[$i]
}
